# feat: add 2022-Q3 data
#
# Before: sheets = [ "总计", "2022-Q1" ]
# After:  sheets = [ "总计", "2022-Q3" (new data), "2022-Q1" (old data, unchanged) ]
#
# Strategy:
#  1. Duplicate the existing "2022-Q1" sheet (sheetId stays on the original,
#     the duplicate gets a brand new sheetId) so the old Q1 fund-holding data
#     is preserved verbatim on its own tab.
#  2. Rename the duplicate back to "2022-Q1".
#  3. Rename the original sheet to "2022-Q3" and replace its data with the
#     new quarter's fund holdings.
#  4. Update the "总计" (totals) sheet: the existing row becomes the 2022-Q3
#     summary row, and a new row is appended below it with the old 2022-Q1
#     summary values.

$wb = $excel.ActiveWorkbook

$xlPasteAll = -4104
$xlPasteFormats = -4122

# ---------------------------------------------------------------------
# 1) Duplicate "2022-Q1" so its current data survives on its own tab.
# ---------------------------------------------------------------------
$q1 = $wb.Worksheets.Item(2)
$q1.Copy($null, $q1)

$q1Copy = $wb.Worksheets.Item(3)

# ---------------------------------------------------------------------
# 2) Turn the original sheet into "2022-Q3" and load its new data. Rename
#    the original FIRST so the duplicate can reclaim the "2022-Q1" name
#    without a transient collision.
# ---------------------------------------------------------------------
$q1.Name = "2022-Q3"
$q1Copy.Name = "2022-Q1"

# Drop the old 3rd data row (new sheet only has 2 data rows) and clear the
# remaining data rows, keeping header row + styles intact.
$q1.Rows.Item(4).Delete()
$q1.Range("A2:H3").ClearContents()

# Header row + A column use the bold/boxed style ("总计" header style) in the
# refreshed sheet.
$wb.Worksheets.Item(1).Range("B1").Copy()
$q1.Range("B1:H1").PasteSpecial($xlPasteFormats)
$wb.Worksheets.Item(1).Range("A2").Copy()
$q1.Range("A2:A3").PasteSpecial($xlPasteFormats)

# Row 2
$q1.Range("A2").Value = 0
$q1.Range("B2:G2").NumberFormat = "@"
$q1.Range("B2").Value = "001628"
$q1.Range("C2").Value = "招商体育文化休闲股票A"
$q1.Range("D2").Value = "2.23"
$q1.Range("E2").Value = "92.42"
$q1.Range("F2").Value = "4.91"
$q1.Range("G2").Value = "0.1095"
$q1.Range("H2").Value = 7

# Row 3
$q1.Range("A3").Value = 1
$q1.Range("B3:G3").NumberFormat = "@"
$q1.Range("B3").Value = "015395"
$q1.Range("C3").Value = "招商体育文化休闲股票C"
$q1.Range("D3").Value = "0.25"
$q1.Range("E3").Value = "92.42"
$q1.Range("F3").Value = "4.91"
$q1.Range("G3").Value = "0.0123"
$q1.Range("H3").Value = 7

# ---------------------------------------------------------------------
# 3) Update the "总计" sheet: push the existing (2022-Q1) row down to row 3,
#    then overwrite row 2 with the new 2022-Q3 summary values.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item(1)

$total.Range("B2:D2").Copy()
$total.Range("B3").PasteSpecial($xlPasteAll)

$total.Range("A2").Copy()
$total.Range("A3").PasteSpecial($xlPasteFormats)
$total.Range("A3").Value = 1

$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 0.12

$excel.CutCopyMode = 0
